# UPDATE ADD _ ANSIBLE _ ROLE
# Append 10 new "issue -> note" rows (A: index 1..10, B: description) to the
# bottom of the running sheet, right after the existing row 18.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arrow = [char]0x2192

$notes = @(
    "percona db $arrow false",
    "waite time $arrow pause",
    "log docker $arrow failse ",
    "Looix rabbitmq voi lenh echo ",
    "haproxy thieu sysctl -p",
    "deployment va service setup index  $arrow tran queue",
    "Thieu mkdir nova",
    "openvswitch and eno2 bus",
    "clean skip",
    "add ansible service "
)

$startRow = 19
for ($i = 0; $i -lt $notes.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $i + 1
    $ws.Cells.Item($row, 2).Value = $notes[$i]
}

# Match the author's final selection/scroll position recorded in the diff.
$null = $ws.Range("B28").Select()
